# Update the "tri proximity" summary tables (Means & Standard Deviations)
# with refreshed statistics, per commit "all updates to include tri proximity tables".

$wb = $excel.ActiveWorkbook

# --- Sheet: Means ---
$ws1 = $wb.Worksheets.Item("Means")

$ws1.Range("E2").Value = 32
$ws1.Range("F2").Value = 51
$ws1.Range("G2").Value = 69

$ws1.Range("E3").Value = 63
$ws1.Range("F3").Value = 37
$ws1.Range("G3").Value = 22

$ws1.Range("E4").Value = 5.4
$ws1.Range("F4").Value = 12
$ws1.Range("G4").Value = 8.9

$ws1.Range("E5").Value = 18
$ws1.Range("F5").Value = 35

$ws1.Range("F7").Value = 13
$ws1.Range("G7").Value = 8.3

$ws1.Range("F8").Value = 11
$ws1.Range("G8").Value = 7.4

$ws1.Range("E9").Value = 43
$ws1.Range("F9").Value = 51

# --- Sheet: Standard Deviations ---
$ws2 = $wb.Worksheets.Item("Standard Deviations")

$ws2.Range("E2").Value = 27
$ws2.Range("F2").Value = 30

$ws2.Range("E3").Value = 30
$ws2.Range("F3").Value = 33
$ws2.Range("G3").Value = 30

$ws2.Range("E4").Value = 4.6
$ws2.Range("F4").Value = 9.2
$ws2.Range("G4").Value = 9.3

$ws2.Range("E5").Value = 20
$ws2.Range("F5").Value = 27
$ws2.Range("G5").Value = 22

$ws2.Range("E6").Value = 15

$ws2.Range("E7").Value = 8.7
$ws2.Range("F7").Value = 9.1

$ws2.Range("E8").Value = 13
$ws2.Range("F8").Value = 10

$ws2.Range("E9").Value = 5.4

$ws2.Range("E10").Value = 0.000000000000000023
$ws2.Range("F10").Value = 0.014
$ws2.Range("G10").Value = 0.0096
